$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B values to reflect the new figures
$ws.Range("B2").Value = 3065.78
$ws.Range("B3").Value = 13832.84
$ws.Range("B4").Value = 23757.889999999996
$ws.Range("B5").Value = 4691.380000000001
$ws.Range("B6").Value = 5350.51
$ws.Range("B7").Value = 5839.5
$ws.Range("B8").Value = 24867.599999999991
$ws.Range("B9").Value = 4571.6399999999985
$ws.Range("B10").Value = 15178.62
$ws.Range("B11").Value = 17782.680000000004

$ws.Range("B17").Value = 1057.2599999999979
$ws.Range("B18").Value = 18008.259999999998
$ws.Range("B19").Value = 2529.9199999999983
$ws.Range("B20").Value = 26850.919999999995
$ws.Range("B21").Value = 11632.539999999997
$ws.Range("B22").Value = 11983.019999999999
$ws.Range("B23").Value = 7144.079999999999
$ws.Range("B24").Value = 17348.29
$ws.Range("B25").Value = 11630.76
$ws.Range("B26").Value = 34543.11
$ws.Range("B27").Value = 3158.3999999999987
$ws.Range("B28").Value = 23787.680000000004
$ws.Range("B29").Value = 16832.479999999996
$ws.Range("B30").Value = 10295.040000000001
$ws.Range("B31").Value = 12815.04

# Update the selected cell in the sheet view
$ws.Range("S6").Select()
